$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.4846292089395843, 1.3781639410256596, -0.3673941227252412, 2.110954202427628)"
$ws.Range("C2").Value = "JSU(-1.4653550008576495, 1.1465709972337499, 1.9928724526249957, 4.162373694393471)"
$ws.Range("D2").Value = "NCT(2.7019086712247686, 1.511770239734573, -0.9275843991507812, 2.4580439015409343)"
$ws.Range("E2").Value = "NIG(1.0445799987677407, 0.8281669264393693, 4.021273944960039, 4.5858575524410545)"
